# "commit ke3 yg ini dr publik dari cmd, isi file berubah"
# Two new cells of data were added below the existing table on Sheet1
# (row 60 was left blank, new entries land on row 61), and the
# selection ends up on C62 afterwards (as if the user tabbed/entered
# past the last filled cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B61").Value = "text baru"
$ws.Range("C61").Value = "utk ujicoba git"

$ws.Range("C62").Select()
